{"js": "// Ajustes no css de priceslider e do modal\n// Colors the two \"to-do\" paragraphs (the \"fechar modal\" bullet and the\n// \"priceslider thumbs\" bullet, including all of its runs) with font color\n// 00B0F0, matching the author's highlight edit.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Unique substrings identifying the two paragraphs touched by the edit.\nconst targets = [\n  \"Criar bot\u00e3o para fechar modal ou simplesmente reposicion\u00e1-lo\",\n  \"Aumentar o tamanho das thumbs do priceslider no filter\"\n];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text;\n  if (targets.some((t) => text.indexOf(t) !== -1)) {\n    // Setting font.color on the paragraph applies it to the paragraph mark\n    // (pPr/rPr) as well as every run in the paragraph.\n    paragraph.font.color = \"#00B0F0\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Ajustes no css de priceslider e do modal\n# Colors the two \"to-do\" paragraphs (the \"fechar modal\" bullet and the\n# \"priceslider thumbs\" bullet, including all of its runs) with font color\n# 00B0F0, matching the author's highlight edit.\n#\n# Note: Word's COM Font.Color is a COLORREF (0x00BBGGRR), the reverse byte\n# order of the OOXML hex (RRGGBB). 0x00B0F0 (RRGGBB) -> 0xF0B000 (BBGGRR).\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if (($t -like \"*Criar bot\u00e3o para fechar modal ou simplesmente reposicion\u00e1-lo*\") -or `\n        ($t -like \"*Aumentar o tamanho das thumbs do priceslider no filter*\")) {\n        # Setting Font.Color on the paragraph's Range colors the paragraph\n        # mark (pPr/rPr) as well as every run within the paragraph.\n        $p.Range.Font.Color = 0xF0B000\n    }\n}\n"}
